$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 'D' = 45007; 'L' = 'Segunda'; 'M' = 160; 'N' = 27000; 'O' = 28000; 'P' = 27500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1375; 'T' = 20 }
    3 = @{ 'D' = 45014; 'L' = 'Segunda'; 'M' = 200; 'N' = 24000; 'O' = 25000; 'P' = 24500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1225; 'T' = 20 }
    4 = @{ 'D' = 44965; 'L' = 'Primera'; 'M' = 100; 'N' = 34000; 'O' = 35000; 'P' = 34600; 'Q' = '$/caja 18 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1922; 'T' = 18 }
    5 = @{ 'D' = 44965; 'L' = 'Segunda'; 'M' = 120; 'N' = 32000; 'O' = 33000; 'P' = 32333; 'Q' = '$/caja 18 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1796; 'T' = 18 }
    6 = @{ 'D' = 44972; 'L' = 'Segunda'; 'M' = 140; 'N' = 27000; 'O' = 28000; 'P' = 27429; 'Q' = '$/caja 18 kilos'; 'R' = 'Región Metropolitana'; 'S' = 1524; 'T' = 18 }
    7 = @{ 'D' = 44993; 'L' = 'Segunda'; 'M' = 130; 'N' = 25000; 'O' = 26000; 'P' = 25462; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1273; 'T' = 20 }
    8 = @{ 'D' = 44643; 'L' = 'Primera'; 'M' = 160; 'N' = 28000; 'O' = 30000; 'P' = 29000; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1450; 'T' = 20 }
    9 = @{ 'D' = 44650; 'L' = 'Primera'; 'M' = 160; 'N' = 31000; 'O' = 32000; 'P' = 31500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1575; 'T' = 20 }
    10 = @{ 'D' = 44650; 'L' = 'Segunda'; 'M' = 250; 'N' = 29000; 'O' = 30000; 'P' = 29500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1475; 'T' = 20 }
    11 = @{ 'D' = 44671; 'L' = 'Segunda'; 'M' = 200; 'N' = 29000; 'O' = 30000; 'P' = 29500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1475; 'T' = 20 }
    12 = @{ 'D' = 44979; 'L' = 'Segunda'; 'M' = 250; 'N' = 29000; 'O' = 30000; 'P' = 29500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1475; 'T' = 20 }
    13 = @{ 'D' = 45028; 'L' = 'Segunda'; 'M' = 200; 'N' = 21000; 'O' = 22000; 'P' = 21500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1075; 'T' = 20 }
    14 = @{ 'D' = 44636; 'L' = 'Primera'; 'M' = 200; 'N' = 29000; 'O' = 30000; 'P' = 29500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1475; 'T' = 20 }
    15 = @{ 'D' = 44664; 'L' = 'Segunda'; 'M' = 150; 'N' = 29000; 'O' = 30000; 'P' = 29500; 'Q' = '$/caja 18 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1639; 'T' = 18 }
    16 = @{ 'D' = 45021; 'L' = 'Segunda'; 'M' = 250; 'N' = 22000; 'O' = 23000; 'P' = 22500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1125; 'T' = 20 }
    17 = @{ 'D' = 44679; 'L' = 'Segunda'; 'M' = 200; 'N' = 29000; 'O' = 30000; 'P' = 29500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1475; 'T' = 20 }
    18 = @{ 'D' = 44679; 'L' = 'Tercera'; 'M' = 200; 'N' = 24000; 'O' = 25000; 'P' = 24500; 'Q' = '$/caja 20 kilos'; 'R' = 'Región de Coquimbo'; 'S' = 1225; 'T' = 20 }
}

$colIndex = @{ D=4; L=12; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20 }

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $colNum = $colIndex[$col]
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$col]
    }
}
